# monitoring_df.xlsx - Sheet1 holds a small "metric name -> value per scenario"
# table (columns B/C/D/E = scenarios 2020/JTMT/iplan/bau). This change inserts a
# "Taz_Palestinian" row (right after Taz_U_Orthodox), inserts "univ_2020" /
# "univ_2050" rows (right after Taz_Jewish), and appends new rows implementing
# add_taz_num_is_unique_row, add_taz_num_count_row & add_percentage_growth_row:
# zonetype_*, jeru_metro_*, in_jerusal_*, taz_num_is_unique, taz_num_count and
# the five "percentage growth" rows.
#
# Rather than shifting rows with native Insert() (which Excel tends to stamp with
# fresh/ad-hoc cell styles), every row below "Taz_U_Orthodox" is simply
# (re)written in place at its final row number - rows 7-13 already existed
# (so they keep their original column-A style), rows 14-30 are brand new and
# have the bold/centered/bordered label style copied onto column A explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row#, isNewRow, A-label, then (value, isTextNumber) for columns B, C, D, E.
# isTextNumber marks values such as "81,184" that must stay text (matching the
# comma-grouped numbers already used elsewhere in the sheet) instead of being
# auto-converted to a number by Excel.
$rowData = @(
     ,(7,  $false, "Taz_Palestinian",                   157,       $false, 157,       $false, 157,       $false, 157,       $false)
     ,(8,  $false, "Taz_arabs_behined_seperation_wall",  5,        $false, 5,         $false, 5,         $false, 5,         $false)
     ,(9,  $false, "Taz_Arab",                           86,       $false, 86,        $false, 86,        $false, 86,        $false)
     ,(10, $false, "Taz_Jewish",                         482,      $false, 482,       $false, 472,       $false, 450,       $false)
     ,(11, $false, "univ_2020",                          "81,184", $true,  "81,184",  $true,  "81,184",  $true,  "81,184",  $true)
     ,(12, $false, "univ_2050",                          "NaN",    $false, "81,184",  $true,  "81,184",  $true,  "81,184",  $true)
     ,(13, $false, "student_2020",                       "412,854",$true,  "412,854", $true,  "412,854", $true,  "412,854", $true)
     ,(14, $true,  "student_2050",                       "NaN",    $false, "903,816", $true,  "833,848", $true,  "963,936", $true)
     ,(15, $true,  "student_yeshiva_and_kollim_2020",     "110,880",$true,  "110,880", $true,  "110,880", $true,  "110,880", $true)
     ,(16, $true,  "student_yeshiva_and_kollim_2050",     "NaN",   $false, "298,512", $true,  "328,614", $true,  "399,184", $true)
     ,(17, $true,  "zonetype_Jerusalem",                  580,     $false, 580,       $false, 580,       $false, 580,       $false)
     ,(18, $true,  "zonetype_Judea and Samaria",          288,     $false, 288,       $false, 288,       $false, 288,       $false)
     ,(19, $true,  "zonetype_Ramla",                      36,      $false, 36,        $false, 36,        $false, 36,        $false)
     ,(20, $true,  "jeru_metro_0",                        74,      $false, 74,        $false, 74,        $false, 74,        $false)
     ,(21, $true,  "jeru_metro_1",                        830,     $false, 830,       $false, 830,       $false, 830,       $false)
     ,(22, $true,  "in_jerusal_yes",                      460,     $false, 460,       $false, 460,       $false, 460,       $false)
     ,(23, $true,  "in_jerusal_no",                       444,     $false, 444,       $false, 444,       $false, 444,       $false)
     ,(24, $true,  "taz_num_is_unique",                   $true,   $false, $true,     $false, $true,     $false, $true,     $false)
     ,(25, $true,  "taz_num_count",                       904,     $false, 904,       $false, 904,       $false, 904,       $false)
     ,(26, $true,  "percentage growth 2020-2025",         "NaN",   $false, 288,       $false, 283,       $false, 285,       $false)
     ,(27, $true,  "percentage growth 2025-2030",         "NaN",   $false, 110,       $false, 109,       $false, 110,       $false)
     ,(28, $true,  "percentage growth 2030-2035",         "NaN",   $false, 109,       $false, 108,       $false, 109,       $false)
     ,(29, $true,  "percentage growth 2035-2040",         "NaN",   $false, 109,       $false, 107,       $false, 109,       $false)
     ,(30, $true,  "percentage growth 2045-2050",         "NaN",   $false, 94,        $false, 93,        $false, 92,        $false)
)

$cols = @("B", "C", "D", "E")

foreach ($r in $rowData) {
    $rowNum   = $r[0]
    $isNewRow = $r[1]
    $label    = $r[2]

    $ws.Range("A$rowNum").Value = $label

    if ($isNewRow) {
        # Brand-new row: pull in the label style (bold font, thin border,
        # center/top alignment) used by every other metric-name cell.
        $ws.Range("A6").Copy()
        $ws.Range("A$rowNum").PasteSpecial(-4122)
    }

    for ($i = 0; $i -lt 4; $i++) {
        $value        = $r[3 + $i * 2]
        $isTextNumber = $r[4 + $i * 2]
        $cell         = $ws.Range("$($cols[$i])$rowNum")

        if ($isTextNumber) {
            # Force text so "81,184" etc. doesn't get auto-converted to a number.
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
